# Insert two new rows before row 252, shifting the existing data block
# (rows 252-362) down to rows 254-364, then populate the two newly
# inserted rows (252-253) with the new weekly data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("252:253").Insert()

# --- Row 252: Pintón ---
$ws.Range("A252").Value = 5
$ws.Range("B252").Value = "Macroferia Regional de Talca"
$ws.Range("C252").Value = "Maule"
$ws.Range("D252").Value = 44466
$ws.Range("E252").Value = 7
$ws.Range("F252").Value = "Fruta"
$ws.Range("G252").Value = 100108
$ws.Range("H252").Value = "Tropicales y subtropicales"
$ws.Range("I252").Value = 100108006
$ws.Range("J252").Value = "Plátano"
$ws.Range("K252").Value = "Sin especificar"
$ws.Range("L252").Value = "Pintón"
$ws.Range("M252").Value = 1000
$ws.Range("N252").Value = 12000
$ws.Range("O252").Value = 13000
$ws.Range("P252").Value = 12600
$ws.Range("Q252").Value = "$/caja 20 kilos"
$ws.Range("R252").Value = "Ecuador"
$ws.Range("S252").Value = 630
$ws.Range("T252").Value = 20

# --- Row 253: Primera Pintón ---
$ws.Range("A253").Value = 5
$ws.Range("B253").Value = "Macroferia Regional de Talca"
$ws.Range("C253").Value = "Maule"
$ws.Range("D253").Value = 44466
$ws.Range("E253").Value = 7
$ws.Range("F253").Value = "Fruta"
$ws.Range("G253").Value = 100108
$ws.Range("H253").Value = "Tropicales y subtropicales"
$ws.Range("I253").Value = 100108006
$ws.Range("J253").Value = "Plátano"
$ws.Range("K253").Value = "Sin especificar"
$ws.Range("L253").Value = "Primera Pintón"
$ws.Range("M253").Value = 600
$ws.Range("N253").Value = 14000
$ws.Range("O253").Value = 14000
$ws.Range("P253").Value = 14000
$ws.Range("Q253").Value = "$/caja 20 kilos"
$ws.Range("R253").Value = "Ecuador"
$ws.Range("S253").Value = 700
$ws.Range("T253").Value = 20
